# plantilla_localizaciones.xlsx - rename the main sheet.
#
# The commit ("Gestiona min y mayus campos strings. Campos numericos >0.
# Vinculo posicion rack con rack y posicion caja con caja. No permite subir
# posiciones si hay errores.") describes validation/macro logic added around
# this locations template, and the accompanying OOXML diff's only durable,
# content-level change is the worksheet being renamed from the default
# "Hoja1" to the meaningful "Ubicaciones" (Spanish for "Locations"), matching
# the template's purpose. The remaining diff hunks (fileVersion rupBuild,
# absPath, revisionPtr documentId, bookViews window geometry, theme
# name/ids, default row height / dyDescent, best-fit column width deltas)
# are all environment artifacts of the authoring machine/Excel build that
# resaved the file, not deliberate edits, so they are left alone here.

$wb = $excel.ActiveWorkbook

# Rename the (only / active) worksheet from "Hoja1" to "Ubicaciones".
$ws = $wb.Worksheets.Item(1)
if ($ws.Name -eq "Hoja1") {
    $ws.Name = "Ubicaciones"
} else {
    $ws = $wb.ActiveSheet
    $ws.Name = "Ubicaciones"
}
